# Insert a new data row into the "Vega Monumental Concepción - Cebollín" log.
# A new row is inserted at row 40 (shifting the existing rows 40-134 down to
# 41-135) and populated with a fresh reading for 2023-08-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything currently on/after row 40 down by one row.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new observation. Columns not
# listed in the diff (A, B, C, E, F, G, H, I, R) carry the same constant
# values used throughout this market/category sheet.
$ws.Range("A40").Value = 11
$ws.Range("B40").Value = "Vega Monumental Concepción"
$ws.Range("C40").Value = "Bíobío"
$ws.Range("D40").Value = 45162
$ws.Range("D40").NumberFormat = $ws.Range("D39").NumberFormat
$ws.Range("E40").Value = 8
$ws.Range("F40").Value = 100112037
$ws.Range("G40").Value = "Cebollín"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Primera"
$ws.Range("J40").Value = 60
$ws.Range("K40").Value = 5000
$ws.Range("L40").Value = 5500
$ws.Range("M40").Value = 5250
$ws.Range("N40").Value = "`$/paquete 36 unidades"
$ws.Range("O40").Value = "Región Metropolitana"
$ws.Range("P40").Value = 146
$ws.Range("Q40").Value = 36
$ws.Range("R40").Value = "Hortaliza"
